$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $cellRef, $text)
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-CellText $ws "D2" "54.851.01"
Set-CellText $ws "E2" "  +5.55%  "
Set-CellText $ws "D3" "2.432.35"
Set-CellText $ws "E3" "  +6.66%  "
Set-CellText $ws "D4" "1.00"
Set-CellText $ws "E4" "  -0.09%  "
Set-CellText $ws "D5" "480.79"
Set-CellText $ws "E5" "  +9.26%  "
Set-CellText $ws "D6" "138.59"
Set-CellText $ws "E6" "  +16.43%  "
Set-CellText $ws "D7" "0.997"
Set-CellText $ws "E7" "  -0.28%  "
Set-CellText $ws "E8" "  +8.56%  "
Set-CellText $ws "D9" "2.454.66"
Set-CellText $ws "E9" "  +7.01%  "
Set-CellText $ws "D10" "0.0967"
Set-CellText $ws "E10" "  +12.53%  "
Set-CellText $ws "E11" "  +3.87%  "
Set-CellText $ws "D12" "0.324"
Set-CellText $ws "E12" "  +8.35%  "
Set-CellText $ws "E13" "  +2.18%  "
Set-CellText $ws "D14" "2.873.28"
Set-CellText $ws "E14" "  +7.58%  "
Set-CellText $ws "D15" "55.030.90"
Set-CellText $ws "E15" "  +5.84%  "
Set-CellText $ws "D16" "20.46"
Set-CellText $ws "E16" "  +10.21%  "
Set-CellText $ws "E17" "  +16.13%  "
Set-CellText $ws "D18" "2.448.21"
Set-CellText $ws "E18" "  +5.51%  "
Set-CellText $ws "D19" "4.33"
Set-CellText $ws "E19" "  +10.76%  "
Set-CellText $ws "D20" "314.23"
Set-CellText $ws "E20" "  +6.72%  "
Set-CellText $ws "D21" "9.82"
Set-CellText $ws "E21" "  +11.91%  "
Set-CellText $ws "D22" "0.997"
Set-CellText $ws "E22" "  -0.14%  "
Set-CellText $ws "D23" "5.62"
Set-CellText $ws "E23" "  +10.46%  "
Set-CellText $ws "D24" "57.17"
Set-CellText $ws "E24" "  +7.59%  "
Set-CellText $ws "E25" "  +14.73%  "
Set-CellText $ws "E26" "  +11.73%  "
Set-CellText $ws "E27" "  -6.22%  "
Set-CellText $ws "D28" "2.541.52"
Set-CellText $ws "E28" "  +5.98%  "
Set-CellText $ws "D29" "7.32"
Set-CellText $ws "E29" "  +7.47%  "
Set-CellText $ws "D30" "0.0₃0774"
Set-CellText $ws "E30" "  +19.77%  "
Set-CellText $ws "D31" "0.999"
Set-CellText $ws "E31" "  +0.09%  "
Set-CellText $ws "D32" "148.57"
Set-CellText $ws "E32" "  +4.89%  "
Set-CellText $ws "D33" "17.93"
Set-CellText $ws "E33" "  +7.93%  "
Set-CellText $ws "E34" "  +11.64%  "
Set-CellText $ws "D35" "5.14"
Set-CellText $ws "E35" "  +10.52%  "
Set-CellText $ws "D36" "1.12"
Set-CellText $ws "E36" "  +13.98%  "
Set-CellText $ws "E37" "  +7.80%  "
Set-CellText $ws "D38" "0.847"
Set-CellText $ws "E38" "  +5.19%  "
Set-CellText $ws "D39" "0.992"
Set-CellText $ws "E39" "  -0.56%  "
Set-CellText $ws "D40" "33.14"
Set-CellText $ws "E40" "  +4.45%  "
Set-CellText $ws "E41" "  +10.85%  "
Set-CellText $ws "D42" "0.0546"
Set-CellText $ws "E42" "  +9.00%  "
Set-CellText $ws "D43" "0.598"
Set-CellText $ws "E43" "  +7.70%  "
Set-CellText $ws "E44" "  +12.05%  "
Set-CellText $ws "D45" "10.10"
Set-CellText $ws "E45" "  -0.51%  "
Set-CellText $ws "D46" "4.64"
Set-CellText $ws "E46" "  +13.63%  "
Set-CellText $ws "D47" "254.31"
Set-CellText $ws "E47" "  +29.33%  "
Set-CellText $ws "D48" "0.0901"
Set-CellText $ws "E48" "  +10.95%  "
Set-CellText $ws "B49" "Maker"
Set-CellText $ws "C49" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-CellText $ws "D49" "1.934.61"
Set-CellText $ws "E49" "  +1.99%  "
Set-CellText $ws "B50" "VeChain"
Set-CellText $ws "C50" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText $ws "D50" "0.0222"
Set-CellText $ws "E50" "  +10.05%  "
Set-CellText $ws "D51" "17.03"
Set-CellText $ws "E51" "  +10.32%  "
